# Apply the changes described by the diff to slide 2 of the presentation.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# EMU -> Points helpers. PowerPoint COM geometry properties (Left/Top/
# Width/Height) are expressed in points (1 pt = 12700 EMU) and this host's
# internal float32 storage truncates (rounds down) on the points -> EMU
# trip, so we bias by +0.5 EMU worth of points to land on the exact target
# EMU value. TextFrame margin properties (MarginTop/Left/Right/Bottom), by
# contrast, round to the nearest EMU, so no bias is needed there.
function EmuToPt($emu) {
    return ($emu + 0.5) / 12700
}
function EmuToPtMargin($emu) {
    return $emu / 12700
}

# --- Shape "CustomShape 10" (background card behind "Correlation coefficient") ---
# ext cy: 6858000 -> 7014726 (offset unchanged)
$shp10 = $s.Shapes.Item("CustomShape 10")
$shp10.Height = EmuToPt 7014726

# --- Shape "CustomShape 16" (rounded header "Categorical response, ...") ---
# off y: 7636376 -> 7772400 (ext unchanged)
$shp16 = $s.Shapes.Item("CustomShape 16")
$shp16.Top = EmuToPt 7772400

# --- Shape "CustomShape 20" (card with "Logistic regression" code sample) ---
# off y: 8395615 -> 8537383 ; ext cy: 1877418 -> 1735649
$shp20 = $s.Shapes.Item("CustomShape 20")
$shp20.Top = EmuToPt 8537383
$shp20.Height = EmuToPt 1735649

# bodyPr tIns: 91440 -> 45720
$shp20.TextFrame.MarginTop = EmuToPtMargin 45720

# First paragraph: "Logistic regression" (single run) becomes two runs,
# "Logistic " + "regression", with the paragraph gaining 150% line spacing.
$tf20 = $shp20.TextFrame
$tr20 = $tf20.TextRange
$para1 = $tr20.Paragraphs(1, 1)
$para1.ParagraphFormat.SpaceWithin = 1.5

$splitPoint = "Logistic ".Length
$firstChunk = $para1.Characters(1, $splitPoint)
$firstChunk.Text = "Logistic "
